$wb = $excel.ActiveWorkbook

# --- Rename "Include from unknown[ N]" sheets to "Include #N" ---
$oldNames = @("Include from unknown", "Include from unknown 2", "Include from unknown 3", "Include from unknown 4", "Include from unknown 5")
for ($i = 0; $i -lt $oldNames.Length; $i++) {
    $wb.Worksheets.Item($oldNames[$i]).Name = "Include #" + $i
}

# --- Update Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Bump published version
$ws.Range("B3").Value = "1.0.1"

# Update contact display text
$ws.Range("B10").Value = "MedCom (http://www.medcom.dk)"

# Insert a new "Jurisdiction" row right after "Contact" (row 10), pushing
# Description/Purpose/Copyright/Immutable down by one row.
$ws.Rows.Item(11).Insert()

# Match the formatting used by the other data rows (copy style from the row above).
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
